# Update cryptos list with latest scraped price/volume data (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '36.932.46'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.67%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.980.67'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.80%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '245.51'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.36%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.631'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.67%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '61.03'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +3.10%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.385'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.04%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0805'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -1.62%  '
$ws.Range('E11').Value = '  +0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '14.79'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +7.42%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.849'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.17%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '22.14'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -1.58%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.270.49'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +0.69%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '5.48'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +3.68%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '1.983.92'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +0.90%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '36.815.05'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.50%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '70.20'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.43%  '
$ws.Range('E20').Value = '  +0.17%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.17'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +1.48%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '230.22'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.43%  '
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('E24').Value = '  +2.68%  '
$ws.Range('E25').Value = '  +0.35%  '
$ws.Range('E26').Value = '  +1.07%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.32'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.11%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '163.87'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.89%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.55'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +0.64%  '
$ws.Range('E30').Value = '  +19.28%  '
$ws.Range('E31').Value = '  +1.59%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.86'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +2.81%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0622'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +0.01%  '
$ws.Range('E34').Value = '  +4.97%  '
$ws.Range('E35').Value = '  -0.34%  '
$ws.Range('E36').Value = '  -0.09%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.36'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -1.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -9.45%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0998'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.42%  '
$ws.Range('E41').Value = '  +0.56%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.18'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.52%  '
$ws.Range('E43').Value = '  +0.14%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '16.31'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +0.44%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.370.42'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '90.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.31%  '
$ws.Range('E47').Value = '  +0.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '7.26'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.16%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '46.43'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +5.27%  '
$ws.Range('B50').Value = 'MXToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.82'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.66%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.98'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +11.42%  '
